$wb = $excel.ActiveWorkbook

# zh-cn sheet
$wsZh = $wb.Worksheets.Item("zh-cn")
$wsZh.Range("E2").Value = "2016-03-22 21:13:51"
$wsZh.Range("H2").Value = "2016-03-22 21:14:19"

# de-de sheet
$wsDe = $wb.Worksheets.Item("de-de")
$wsDe.Range("E2").Value = "2016-03-22 21:13:55"
$wsDe.Range("H2").Value = "2016-03-22 21:14:26"
